$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("D7").Value = "MEC-2B-Máquinas Térmicas e de Fluxo"
$ws.Range("D8").Value = "MEC-2B-Máquinas Térmicas e de Fluxo"
